$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-22 20:48:18"
$ws.Range("K2").Value = "12.9 MJ/m2"
$ws.Range("E3").Value = "2026-02-22 20:48:21"
$ws.Range("O3").Value = "4.5 °C"
$ws.Range("E4").Value = "2026-02-22 20:48:23"
$ws.Range("H4").Value = "'64%"
$ws.Range("O4").Value = "12.4 °C"
$ws.Range("E5").Value = "2026-02-22 20:48:26"
$ws.Range("H5").Value = "'29%"
$ws.Range("O5").Value = "6.1 °C"
$ws.Range("E6").Value = "2026-02-22 20:48:28"
$ws.Range("E7").Value = "2026-02-22 20:48:31"
$ws.Range("E8").Value = "2026-02-22 20:48:33"
$ws.Range("E9").Value = "2026-02-22 20:48:36"
$ws.Range("H9").Value = "'76%"
$ws.Range("O9").Value = "11.0 °C"
$ws.Range("E10").Value = "2026-02-22 20:48:37"
$ws.Range("O10").Value = "10.1 °C"
$ws.Range("E11").Value = "2026-02-22 20:48:38"
$ws.Range("E12").Value = "2026-02-22 20:48:39"
$ws.Range("O12").Value = "9.9 °C"
$ws.Range("E13").Value = "2026-02-22 20:48:40"
$ws.Range("E14").Value = "2026-02-22 20:48:41"
$ws.Range("H14").Value = "'73%"
$ws.Range("O14").Value = "12.0 °C"
$ws.Range("E15").Value = "2026-02-22 20:48:42"
$ws.Range("O15").Value = "10.9 °C"
$ws.Range("E16").Value = "2026-02-22 20:48:43"
$ws.Range("E17").Value = "2026-02-22 20:48:44"
$ws.Range("E18").Value = "2026-02-22 20:48:46"
$ws.Range("J18").Value = "1027.5 hPa"
$ws.Range("E19").Value = "2026-02-22 20:48:47"
$ws.Range("O19").Value = "12.2 °C"
$ws.Range("E20").Value = "2026-02-22 20:48:49"
$ws.Range("K20").Value = "16.2 MJ/m2"
$ws.Range("O20").Value = "4.1 °C"
$ws.Range("E21").Value = "2026-02-22 20:48:51"
$ws.Range("K21").Value = "15.4 MJ/m2"
$ws.Range("E22").Value = "2026-02-22 20:48:54"
$ws.Range("E23").Value = "2026-02-22 20:48:56"
$ws.Range("E24").Value = "2026-02-22 20:48:59"
$ws.Range("H24").Value = "'83%"
$ws.Range("E25").Value = "2026-02-22 20:49:01"
$ws.Range("O25").Value = "7.1 °C"
$ws.Range("E26").Value = "2026-02-22 20:49:04"
$ws.Range("J26").Value = "1026.1 hPa"
$ws.Range("O26").Value = "11.3 °C"
$ws.Range("E27").Value = "2026-02-22 20:49:06"
$ws.Range("E28").Value = "2026-02-22 20:49:09"
$ws.Range("H28").Value = "'64%"
$ws.Range("O28").Value = "10.5 °C"
$ws.Range("E29").Value = "2026-02-22 20:49:11"
$ws.Range("O29").Value = "9.9 °C"
$ws.Range("E30").Value = "2026-02-22 20:49:14"
$ws.Range("H30").Value = "'73%"
$ws.Range("E31").Value = "2026-02-22 20:49:16"
$ws.Range("O31").Value = "14.7 °C"
$ws.Range("E32").Value = "2026-02-22 20:49:19"
$ws.Range("H32").Value = "'70%"
$ws.Range("O32").Value = "6.1 °C"
$ws.Range("E33").Value = "2026-02-22 20:49:21"
$ws.Range("E34").Value = "2026-02-22 20:49:24"
$ws.Range("E35").Value = "2026-02-22 20:49:26"
$ws.Range("E36").Value = "2026-02-22 20:49:29"
$ws.Range("E37").Value = "2026-02-22 20:49:31"
$ws.Range("E38").Value = "2026-02-22 20:49:34"
$ws.Range("E39").Value = "2026-02-22 20:49:36"
$ws.Range("H39").Value = "'25%"
$ws.Range("L39").Value = "27.7 km/h - 314º 20:10 TU"
$ws.Range("M39").Value = "7.3 °C 20:15 TU"
$ws.Range("E40").Value = "2026-02-22 20:49:39"
$ws.Range("O40").Value = "10.1 °C"
$ws.Range("E41").Value = "2026-02-22 20:49:41"
$ws.Range("H41").Value = "'77%"
$ws.Range("J41").Value = "1027.6 hPa"
$ws.Range("K41").Value = "15.5 MJ/m2"
$ws.Range("E42").Value = "2026-02-22 20:49:44"
$ws.Range("H42").Value = "'81%"
$ws.Range("O42").Value = "10.7 °C"
$ws.Range("E43").Value = "2026-02-22 20:49:46"
$ws.Range("E44").Value = "2026-02-22 20:49:48"
$ws.Range("E45").Value = "2026-02-22 20:49:50"
$ws.Range("J45").Value = "1029.0 hPa"
$ws.Range("O45").Value = "8.8 °C"
$ws.Range("E46").Value = "2026-02-22 20:49:53"
